$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the BOM row for U1: part number and footprint changed
$ws.Range("A7").Value = "SN74LVC1G07DCKR"
$ws.Range("C7").Value = "SC-70-5"
$ws.Range("D7").Value = "C7830"

# Update the last-used selection to reflect where the user left off editing
$ws.Range("D8").Select()
